# Daily attendance processing - 2026-02-06 22:58:53 UTC
# Update "Recorded By" (column G) entries from "Administrator, Miss Dina Nasr"
# to "Miss Dina Nasr, Administrator" for the relevant rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(2,4,5,8,17,18,19,20,21,22,23,24,26,28,30,31,34,43,44,45,46,47,48,49,50,52,54,56,57,60,69,70,71,72,73,74,75,76,78,80,81,82,93,94,96,99,101,106,107,108,119,120,122,125,127,132,133,134,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Administrator, Miss Dina Nasr") {
        $cell.Value2 = "Miss Dina Nasr, Administrator"
    }
}
